$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Specification")

$ws.Range("B21").Value = "Telephone number and email address of the applicant."
$ws.Range("B25").Value = "Name and contact information for the parties making the application."
$ws.Range("B31").Value = "Name and contact information if an agent is being used."
$ws.Range("B35").Value = "Name and contact information if an agent is being used."
$ws.Range("B43").Value = "How any natural habitats on the development site will be improved by the proposed works."
$ws.Range("B58").Value = "Checking whether all the requirements of the form have been met, such as proof of payment or supporting documentation."
$ws.Range("B59").Value = "What community consultation activities have taken place as part of the application"
$ws.Range("B61").Value = "Details of any conflict of interest that may exist between the applicant and planning authority."
$ws.Range("B64").Value = "Signed and dated verification of the application's accuracy."
$ws.Range("B67").Value = "Why demolition is necessary at the development site"
$ws.Range("B68").Value = "Who will be affected by the proposal and whether they have been notified, such as agricultural tenants"
$ws.Range("B83").Value = "Details of pre-application advice received from the planning authority"
$ws.Range("B88").Value = "What development, works or change of use is proposed"
$ws.Range("B93").Value = "Details of any other development proposals made for the site"
$ws.Range("B97").Value = "Where the proposed development will be built."
$ws.Range("B106").Value = "Information to help the planning authority arrange a site visit"
